$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("question_answers")
$ws2 = $wb.Worksheets.Item("outputs")

$c = $ws1.Range("B4")
$c.NumberFormat = "@"
$c.Value = "4"
$c.Style = "Normal"
$c = $ws1.Range("B5")
$c.NumberFormat = "@"
$c.Value = "3"
$c.Style = "Normal"
$c = $ws1.Range("B7")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$c = $ws1.Range("B8")
$c.NumberFormat = "@"
$c.Value = "4"
$c.Style = "Normal"
$c = $ws1.Range("B9")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$c = $ws1.Range("B10")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Style = "Normal"
$c = $ws1.Range("B11")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$c = $ws1.Range("B12")
$c.NumberFormat = "@"
$c.Value = "4"
$c.Style = "Normal"
$c = $ws1.Range("B13")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Style = "Normal"
$c = $ws1.Range("B16")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Style = "Normal"
$c = $ws1.Range("B18")
$c.NumberFormat = "@"
$c.Value = "3"
$c.Style = "Normal"
$c = $ws1.Range("B19")
$c.NumberFormat = "@"
$c.Value = "5"
$c.Style = "Normal"
$c = $ws1.Range("B20")
$c.NumberFormat = "@"
$c.Value = "5"
$c.Style = "Normal"
$c = $ws1.Range("B21")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Style = "Normal"
$c = $ws1.Range("B22")
$c.NumberFormat = "@"
$c.Value = "3"
$c.Style = "Normal"
$c = $ws1.Range("B23")
$c.NumberFormat = "@"
$c.Value = "4"
$c.Style = "Normal"
$c = $ws1.Range("B24")
$c.NumberFormat = "@"
$c.Value = "5"
$c.Style = "Normal"
$c = $ws1.Range("B26")
$c.NumberFormat = "@"
$c.Value = "3"
$c.Style = "Normal"
$c = $ws1.Range("B27")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$c = $ws1.Range("B28")
$c.NumberFormat = "@"
$c.Value = "3"
$c.Style = "Normal"
$c = $ws1.Range("B29")
$c.NumberFormat = "@"
$c.Value = "4"
$c.Style = "Normal"
$c = $ws1.Range("B30")
$c.NumberFormat = "@"
$c.Value = "3"
$c.Style = "Normal"
$c = $ws1.Range("B32")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Style = "Normal"
$c = $ws1.Range("B33")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$c = $ws1.Range("B34")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Style = "Normal"
$c = $ws1.Range("B35")
$c.NumberFormat = "@"
$c.Value = "3"
$c.Style = "Normal"
$c = $ws1.Range("B36")
$c.NumberFormat = "@"
$c.Value = "3"
$c.Style = "Normal"
$c = $ws1.Range("B37")
$c.NumberFormat = "@"
$c.Value = "4"
$c.Style = "Normal"
$c = $ws1.Range("B38")
$c.NumberFormat = "@"
$c.Value = "5"
$c.Style = "Normal"
$c = $ws1.Range("B39")
$c.NumberFormat = "@"
$c.Value = "5"
$c.Style = "Normal"
$c = $ws1.Range("B40")
$c.NumberFormat = "@"
$c.Value = "5"
$c.Style = "Normal"
$c = $ws1.Range("B41")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Style = "Normal"
$c = $ws1.Range("B42")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Style = "Normal"
$c = $ws1.Range("B43")
$c.NumberFormat = "@"
$c.Value = "3"
$c.Style = "Normal"
$c = $ws1.Range("B45")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$c = $ws1.Range("B46")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$c = $ws1.Range("B47")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Style = "Normal"
$c = $ws1.Range("B48")
$c.NumberFormat = "@"
$c.Value = "4"
$c.Style = "Normal"
$c = $ws1.Range("B49")
$c.NumberFormat = "@"
$c.Value = "4"
$c.Style = "Normal"
$c = $ws1.Range("B50")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Style = "Normal"
$c = $ws1.Range("B51")
$c.NumberFormat = "@"
$c.Value = "4"
$c.Style = "Normal"
$c = $ws1.Range("B52")
$c.NumberFormat = "@"
$c.Value = "3"
$c.Style = "Normal"
$c = $ws1.Range("B53")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$c = $ws1.Range("B54")
$c.NumberFormat = "@"
$c.Value = "5"
$c.Style = "Normal"
$c = $ws1.Range("B55")
$c.NumberFormat = "@"
$c.Value = "4"
$c.Style = "Normal"
$c = $ws1.Range("B56")
$c.NumberFormat = "@"
$c.Value = "3"
$c.Style = "Normal"
$c = $ws1.Range("B58")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$c = $ws1.Range("B59")
$c.NumberFormat = "@"
$c.Value = "5"
$c.Style = "Normal"
$c = $ws1.Range("B60")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Style = "Normal"
$c = $ws1.Range("B61")
$c.NumberFormat = "@"
$c.Value = "5"
$c.Style = "Normal"
$c = $ws1.Range("B62")
$c.NumberFormat = "@"
$c.Value = "4"
$c.Style = "Normal"
$c = $ws1.Range("B64")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$c = $ws1.Range("B65")
$c.NumberFormat = "@"
$c.Value = "4"
$c.Style = "Normal"
$c = $ws1.Range("B66")
$c.NumberFormat = "@"
$c.Value = "3"
$c.Style = "Normal"
$c = $ws1.Range("B67")
$c.NumberFormat = "@"
$c.Value = "4"
$c.Style = "Normal"
$c = $ws1.Range("B68")
$c.NumberFormat = "@"
$c.Value = "3"
$c.Style = "Normal"
$c = $ws1.Range("B69")
$c.NumberFormat = "@"
$c.Value = "3"
$c.Style = "Normal"
$c = $ws1.Range("B70")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Style = "Normal"
$c = $ws1.Range("B72")
$c.NumberFormat = "@"
$c.Value = "4"
$c.Style = "Normal"
$c = $ws1.Range("B73")
$c.NumberFormat = "@"
$c.Value = "4"
$c.Style = "Normal"
$c = $ws1.Range("B74")
$c.NumberFormat = "@"
$c.Value = "1"
$c.Style = "Normal"
$c = $ws1.Range("B76")
$c.NumberFormat = "@"
$c.Value = "3"
$c.Style = "Normal"
$c = $ws1.Range("B78")
$c.NumberFormat = "@"
$c.Value = "5"
$c.Style = "Normal"
$c = $ws1.Range("B79")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$c = $ws1.Range("B80")
$c.NumberFormat = "@"
$c.Value = "4"
$c.Style = "Normal"
$c = $ws1.Range("B81")
$c.NumberFormat = "@"
$c.Value = "5"
$c.Style = "Normal"

$ws2.Range("B2").Value = 18
$ws2.Range("B3").Value = 33
$ws2.Range("B4").Value = 32
$ws2.Range("B5").Value = 29
$ws2.Range("B6").Value = 30
$ws2.Range("B7").Value = 25
$ws2.Range("B8").Value = 33
$ws2.Range("B9").Value = 34
